$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 1766.6666
$ws.Range("I4").Value = 1800
$ws.Range("K4").Value = 1800
$ws.Range("M4").Value = -1686
$ws.Range("H18").Value = 9226.941000000001
$ws.Range("I18").Value = 5256.5835
$ws.Range("K18").Value = 5256.5835
$ws.Range("M18").Value = -4972.5835
$ws.Range("H62").Value = 2072.4
$ws.Range("I62").Value = 2119.25
$ws.Range("J62").Value = 1885
$ws.Range("K62").Value = 2119.25
$ws.Range("L62").Value = 1885
$ws.Range("M62").Value = -1495.25
$ws.Range("N62").Value = -3133
$ws.Range("H65").Value = 2072.4
$ws.Range("I65").Value = 2119.25
$ws.Range("J65").Value = 1885
$ws.Range("K65").Value = 10596.25
$ws.Range("L65").Value = 9425
$ws.Range("M65").Value = -7476.25
$ws.Range("N65").Value = -15665
$ws.Range("H76").Value = 2799.6667
$ws.Range("I76").Value = 1700
$ws.Range("J76").Value = 3349.5
$ws.Range("K76").Value = 1700
$ws.Range("L76").Value = 3349.5
$ws.Range("N76").Value = -3979.5
$ws.Range("M76").Value = -1385
$ws.Range("H79").Value = 2799.6667
$ws.Range("I79").Value = 1700
$ws.Range("J79").Value = 3349.5
$ws.Range("K79").Value = 1700
$ws.Range("L79").Value = 3349.5
$ws.Range("N79").Value = -5533.5
$ws.Range("M79").Value = -608
$ws.Range("H98").Value = 1080.375
$ws.Range("I98").Value = 808.5909
$ws.Range("K98").Value = 808.5909
$ws.Range("M98").Value = 689.4091
$ws.Range("H107").Value = 894.6667
$ws.Range("I107").Value = 721.4545000000001
$ws.Range("K107").Value = 721.4545000000001
$ws.Range("M107").Value = 1198.5455
$ws.Range("H113").Value = 39225.625
$ws.Range("I113").Value = 61001
$ws.Range("J113").Value = 2933.3333
$ws.Range("K113").Value = 61001
$ws.Range("L113").Value = 2933.3333
$ws.Range("M113").Value = -57747
$ws.Range("N113").Value = -9441.3333
$ws.Range("H122").Value = 1080.375
$ws.Range("I122").Value = 808.5909
$ws.Range("K122").Value = 2425.7727
$ws.Range("M122").Value = 24.22730000000001
$ws.Range("H132").Value = 1491.2
$ws.Range("I132").Value = 1516.6666
$ws.Range("J132").Value = 1453
$ws.Range("K132").Value = 4549.9998
$ws.Range("L132").Value = 4359
$ws.Range("M132").Value = -2019.9998
$ws.Range("N132").Value = -9419
$ws.Range("H141").Value = 7610.25
$ws.Range("I141").Value = 2328.3333
$ws.Range("K141").Value = 6984.999899999999
$ws.Range("M141").Value = -1804.999899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3765.5352
$ws.Range("I32").Value = 2638.5156
$ws.Range("K32").Value = 2638.5156
$ws.Range("M32").Value = -2351.5156
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("M53").ClearContents()
$ws.Range("N53").ClearContents()
$ws.Range("H110").Value = 128.625
$ws.Range("I110").Value = 132.71428
$ws.Range("K110").Value = 132.71428
$ws.Range("M110").Value = 1912.28572
$ws.Range("H122").Value = 960.86365
$ws.Range("I122").Value = 796.8946999999999
$ws.Range("J122").Value = 1999.3334
$ws.Range("K122").Value = 2390.6841
$ws.Range("L122").Value = 5998.0002
$ws.Range("M122").Value = 59.31590000000006
$ws.Range("N122").Value = -10898.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4385.121
$ws.Range("I134").Value = 4814.857
$ws.Range("J134").Value = 1978.6
$ws.Range("K134").Value = 14444.571
$ws.Range("L134").Value = 5935.799999999999
$ws.Range("M134").Value = -11909.571
$ws.Range("N134").Value = -11005.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1134.5454
$ws.Range("I94").Value = 1047.875
$ws.Range("K94").Value = 1047.875
$ws.Range("M94").Value = -596.875
$ws.Range("H99").Value = 3212.3333
$ws.Range("I99").Value = 2182.2
$ws.Range("K99").Value = 2182.2
$ws.Range("M99").Value = -684.1999999999998
$ws.Range("H105").Value = 1092.5834
$ws.Range("I105").Value = 1110.875
$ws.Range("K105").Value = 1110.875
$ws.Range("M105").Value = 636.125
$ws.Range("H126").Value = 3212.3333
$ws.Range("I126").Value = 2182.2
$ws.Range("K126").Value = 6546.599999999999
$ws.Range("M126").Value = -4076.599999999999
$ws.Range("H132").Value = 1972
$ws.Range("I132").Value = 1508.8889
$ws.Range("K132").Value = 4526.6667
$ws.Range("M132").Value = -1996.6667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 17369.182
$ws.Range("J131").Value = 18159.857
$ws.Range("L131").Value = 54479.571
$ws.Range("N131").Value = -64559.571

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1138.1428
$ws.Range("I113").Value = 1118
$ws.Range("J113").Value = 1153.25
$ws.Range("K113").Value = 1118
$ws.Range("L113").Value = 1153.25
$ws.Range("M113").Value = 1052
$ws.Range("N113").Value = -5493.25
$ws.Range("H132").Value = 2264652.5
$ws.Range("I132").Value = 2960368.8
$ws.Range("J132").Value = 3574.25
$ws.Range("K132").Value = 8881106.399999999
$ws.Range("L132").Value = 10722.75
$ws.Range("M132").Value = -8878576.399999999
$ws.Range("N132").Value = -15782.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2284.7693
$ws.Range("I7").Value = 2100.1667
$ws.Range("J7").Value = 4500
$ws.Range("K7").Value = 2100.1667
$ws.Range("L7").Value = 4500
$ws.Range("M7").Value = -1988.1667
$ws.Range("N7").Value = -4724
$ws.Range("H55").Value = 423.72726
$ws.Range("I55").Value = 352.53333
$ws.Range("J55").Value = 576.2857
$ws.Range("K55").Value = 352.53333
$ws.Range("L55").Value = 576.2857
$ws.Range("M55").Value = -179.53333
$ws.Range("N55").Value = -922.2857
$ws.Range("H61").Value = 2584.7856
$ws.Range("I61").Value = 2231.889
$ws.Range("K61").Value = 2231.889
$ws.Range("M61").Value = -2029.889
$ws.Range("H113").Value = 2584.7856
$ws.Range("I113").Value = 2231.889
$ws.Range("K113").Value = 2231.889
$ws.Range("M113").Value = -61.88900000000012
$ws.Range("H126").Value = 2284.7693
$ws.Range("I126").Value = 2100.1667
$ws.Range("J126").Value = 4500
$ws.Range("K126").Value = 6300.500100000001
$ws.Range("L126").Value = 13500
$ws.Range("M126").Value = -3830.500100000001
$ws.Range("N126").Value = -18440
$ws.Range("H136").Value = 4882.0435
$ws.Range("I136").Value = 3973.8948
$ws.Range("J136").Value = 9195.75
$ws.Range("K136").Value = 11921.6844
$ws.Range("L136").Value = 27587.25
$ws.Range("M136").Value = -9371.6844
$ws.Range("N136").Value = -32687.25
$ws.Range("H138").Value = 50000
$ws.Range("J138").Value = 50000
$ws.Range("L138").Value = 50000
$ws.Range("N138").Value = -60280

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 736.3684
$ws.Range("I107").Value = 418.69232
$ws.Range("J107").Value = 1424.6666
$ws.Range("K107").Value = 1256.07696
$ws.Range("L107").Value = 4273.9998
$ws.Range("M107").Value = 663.9230400000001
$ws.Range("N107").Value = -8113.9998
$ws.Range("H122").Value = 57053.715
$ws.Range("I122").Value = 87662.22
$ws.Range("J122").Value = 1958.4
$ws.Range("K122").Value = 262986.66
$ws.Range("L122").Value = 5875.200000000001
$ws.Range("M122").Value = -260536.66
$ws.Range("N122").Value = -10775.2
$ws.Range("H126").Value = 1897.7142
$ws.Range("I126").Value = 2097.4546
$ws.Range("J126").Value = 1165.3334
$ws.Range("K126").Value = 6292.3638
$ws.Range("L126").Value = 3496.0002
$ws.Range("M126").Value = -3822.3638
$ws.Range("N126").Value = -8436.0002
